$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Placement")

# Establish new shared strings in the same order they first appear in the
# target workbook: "mm" (written via C24 first), then "inches" (via C25),
# then "Mousebites" (via B23) last.

# ---- mm / inches pair #1 (0.5 mm) ----
$ws.Range("C24").Value = "mm"
$ws.Range("C25").Value = "inches"

# ---- Mousebites section header ----
$ws.Range("B23").Value = "Mousebites"

# Row 24: 0.5 mm value
$ws.Range("B24").Value = 0.5

# Row 25: conversion to inches
$ws.Range("B25").Formula = "=B24/25.4"
$ws.Range("B25").NumberFormat = "0.000"

# ---- mm / inches pair #2 (1 mm) ----
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "mm"

$ws.Range("B28").Formula = "=B27/25.4"
$ws.Range("B28").NumberFormat = "0.000"
$ws.Range("C28").Value = "inches"

# ---- mm / inches pair #3 (90 mm) ----
$ws.Range("B31").Value = 90
$ws.Range("C31").Value = "mm"

$ws.Range("B32").Formula = "=B31/25.4"
$ws.Range("B32").NumberFormat = "0.000"
$ws.Range("C32").Value = "inches"

# ---- mm / inches pair #4 (40 mm) ----
$ws.Range("B34").Value = 40
$ws.Range("C34").Value = "mm"

$ws.Range("B35").Formula = "=B34/25.4"
$ws.Range("B35").NumberFormat = "0.000"
$ws.Range("C35").Value = "inches"

# Final selection / active cell as left by the author
$ws.Range("B33").Select()
